# Regenerate the "K" column (column G) values in the save_data sheet.
# These values are produced upstream (Strike# -> K recalculation, std/mean
# regen, s_vals calc) and are simply written back into the worksheet here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K"), for rows 2 through 65 (data rows, row 1 is header).
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 3
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 1
    39 = 1
    40 = 3
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 2
    46 = 0
    47 = 1
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 1
    54 = 0
    55 = 3
    56 = 2
    57 = 2
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 2
    63 = 1
    64 = 0
    65 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
